$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update OS value from "Windows 10" to "Ubuntu"
$ws.Range("D3").Value = "Ubuntu"

# Update Price value from 300 to 250
$ws.Range("D4").Value = 250

# Update the active selection to D4
$ws.Range("D4").Select()

# Update default column width (affects sheetFormatPr defaultColWidth)
$ws.StandardWidth = 11.66015625
